$p = $ppt.ActivePresentation
try {
  $p.ApplyTheme()
  Write-Output "NOARG-OK"
} catch {
  Write-Output "NOARG-ERR: $_"
}
